# Revert "Changed Date on Proposal": the cover page publish date moves
# back from May 15 to May 10 (the bound Year field underneath it stays
# "2013", only the "MMMM d" formatted field's text changes).
#
# The "May 15" text lives inside the cover-page table, which is nested
# inside content controls (w:sdt). In that region Find/Execute reports
# no match (it does not walk into the nested table/content-control),
# so first try the normal Find/Replace path for robustness, then fall
# back to locating the run by scanning Document.Range offsets directly
# (Range *does* reach inside the table/content-control) and assigning
# Range.Text to rewrite it.

$d = $word.ActiveDocument

$oldDateText = "May 15"
$newDateText = "May 10"

$replaced = $false

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = $oldDateText
$find.Replacement.ClearFormatting()
$find.Replacement.Text = $newDateText
$found = $find.Execute($oldDateText, $false, $false, $false, $false, $false, $true, 1, $false, $newDateText, 2)
if ($found) {
    $replaced = $true
}

if (-not $replaced) {
    # Document.Range() uses its own character coordinate space, which
    # does not line up 1:1 with Content.Text.Length (hidden
    # content-control / field boundary marks inflate the .Text string
    # but are not addressable Range offsets), so probe defensively and
    # stop at the first failure instead of trusting a precomputed
    # upper bound.
    function Find-TextOffset($doc, $needle, $scanLimit) {
        $nlen = $needle.Length
        for ($i = 0; $i -lt $scanLimit; $i++) {
            try {
                $candidate = $doc.Range($i, $i + $nlen).Text
            } catch {
                break
            }
            if ($candidate -eq $needle) {
                return $i
            }
        }
        return -1
    }

    $docLen = $d.Content.Text.Length
    $pos = Find-TextOffset $d $oldDateText $docLen
    if ($pos -ge 0) {
        $rng = $d.Range($pos, $pos + $oldDateText.Length)
        $rng.Text = $newDateText
        $replaced = $true
    }
}
